$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.216.06"
$ws.Range("E2").Value = "  +2.48%  "
$ws.Range("D3").Value = "1.589.38"
$ws.Range("E3").Value = "  +1.01%  "
$ws.Range("E4").Value = "  +0.84%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.78%  "
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("E7").Value = "  +0.91%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.99"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.35%  "
$ws.Range("E9").Value = "  -0.22%  "
$ws.Range("E10").Value = "  +0.40%  "
$ws.Range("E11").Value = "  +2.41%  "
$ws.Range("D12").Value = "1.817.81"
$ws.Range("E12").Value = "  +1.12%  "
$ws.Range("D13").Value = "1.593.32"
$ws.Range("E13").Value = "  +1.22%  "
$ws.Range("E14").Value = "  +1.05%  "
$ws.Range("E15").Value = "  -1.01%  "
$ws.Range("D16").Value = "28.265.98"
$ws.Range("E16").Value = "  +2.69%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.24"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.03%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "227.35"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.16%  "
$ws.Range("E19").Value = "  +0.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.03%  "
$ws.Range("E21").Value = "  +0.70%  "
$ws.Range("E22").Value = "  -2.06%  "
$ws.Range("E23").Value = "  -1.63%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.89"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.69%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.20"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("E27").Value = "  -0.46%  "
$ws.Range("E28").Value = "  -1.48%  "
$ws.Range("E29").Value = "  +0.83%  "
$ws.Range("E30").Value = "  -0.86%  "
$ws.Range("E31").Value = "  -0.19%  "
$ws.Range("E32").Value = "  -0.35%  "
$ws.Range("E33").Value = "  -1.56%  "
$ws.Range("D34").Value = "1.398.65"
$ws.Range("E34").Value = "  -4.29%  "
$ws.Range("E35").Value = "  -2.90%  "
$ws.Range("E36").Value = "  -8.13%  "
$ws.Range("E37").Value = "  +1.36%  "
$ws.Range("E38").Value = "  -0.51%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.55"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.61%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.541"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.36%  "
$ws.Range("E41").Value = "  -0.79%  "
$ws.Range("E42").Value = "  +0.73%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.87"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.83%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.59"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.32%  "
$ws.Range("E45").Value = "  +1.22%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "64.28"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.06%  "
$ws.Range("D47").Value = "1.727.26"
$ws.Range("E47").Value = "  +0.99%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.33"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.70%  "
$ws.Range("E49").Value = "  +1.62%  "
$ws.Range("E50").Value = "  +6.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0523"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.44%  "
